# Laptop_build.xlsx -- "Add numpad power draw" edit
#
# The numpad part (row 16 of the "Parts" table) is missing its power draw
# (Amps) figure. Fill in G16 = 0.25 A so the dependent "Watt total" formula
# in I16 (and the SUBTOTAL() rollups in row 26) recalculate correctly.
# Also bump both sheets' zoom back up to 100% and leave the cursor sitting
# on the newly-computed I16 cell.

$wb = $excel.ActiveWorkbook

$wsParts = $wb.Worksheets.Item("Parts")
$wsCooling = $wb.Worksheets.Item("Cooling")

# Numpad draws 0.25A -- add it to the table; Watt total (I16) and the
# SUBTOTAL()s on row 26 recalc automatically from this.
$wsParts.Range("G16").Value = 0.25

# Restore the zoom level on both sheets to 100% (was left at 50%), and
# park the selection on the cell we just finished computing.
$wsParts.Activate()
$wsParts.Range("I16").Select()
$excel.ActiveWindow.Zoom = 100

$wsCooling.Activate()
$excel.ActiveWindow.Zoom = 100

# Leave "Parts" as the active/selected sheet, matching tabSelected="true".
$wsParts.Activate()
